$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the two new columns I0/IF, matching the style of the
# existing header row (bold font, borders, centered) by copying format
# from H1 then re-applying the text values.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I2:J83 (I0 and IF columns hold identical values per row)
$data = New-Object 'object[,]' 82,2
$data[0,0] = 7
$data[0,1] = 7
$data[1,0] = 9
$data[1,1] = 9
$data[2,0] = 8
$data[2,1] = 8
$data[3,0] = 8
$data[3,1] = 8
$data[4,0] = 5
$data[4,1] = 5
$data[5,0] = 9
$data[5,1] = 9
$data[6,0] = 8
$data[6,1] = 8
$data[7,0] = 9
$data[7,1] = 9
$data[8,0] = 8
$data[8,1] = 8
$data[9,0] = 8
$data[9,1] = 8
$data[10,0] = 9
$data[10,1] = 9
$data[11,0] = 9
$data[11,1] = 9
$data[12,0] = 8
$data[12,1] = 8
$data[13,0] = 7
$data[13,1] = 7
$data[14,0] = 7
$data[14,1] = 7
$data[15,0] = 7
$data[15,1] = 7
$data[16,0] = 8
$data[16,1] = 8
$data[17,0] = 9
$data[17,1] = 9
$data[18,0] = 9
$data[18,1] = 9
$data[19,0] = 8
$data[19,1] = 8
$data[20,0] = 7
$data[20,1] = 7
$data[21,0] = 7
$data[21,1] = 7
$data[22,0] = 7
$data[22,1] = 7
$data[23,0] = 7
$data[23,1] = 7
$data[24,0] = 7
$data[24,1] = 7
$data[25,0] = 8
$data[25,1] = 8
$data[26,0] = 8
$data[26,1] = 8
$data[27,0] = 8
$data[27,1] = 8
$data[28,0] = 8
$data[28,1] = 8
$data[29,0] = 9
$data[29,1] = 9
$data[30,0] = 9
$data[30,1] = 9
$data[31,0] = 11
$data[31,1] = 11
$data[32,0] = 8
$data[32,1] = 8
$data[33,0] = 7
$data[33,1] = 7
$data[34,0] = 7
$data[34,1] = 7
$data[35,0] = 8
$data[35,1] = 8
$data[36,0] = 7
$data[36,1] = 7
$data[37,0] = 7
$data[37,1] = 7
$data[38,0] = 8
$data[38,1] = 8
$data[39,0] = 9
$data[39,1] = 9
$data[40,0] = 8
$data[40,1] = 8
$data[41,0] = 7
$data[41,1] = 7
$data[42,0] = 9
$data[42,1] = 9
$data[43,0] = 8
$data[43,1] = 8
$data[44,0] = 9
$data[44,1] = 9
$data[45,0] = 7
$data[45,1] = 7
$data[46,0] = 8
$data[46,1] = 8
$data[47,0] = 6
$data[47,1] = 6
$data[48,0] = 9
$data[48,1] = 9
$data[49,0] = 9
$data[49,1] = 9
$data[50,0] = 8
$data[50,1] = 8
$data[51,0] = 6
$data[51,1] = 6
$data[52,0] = 8
$data[52,1] = 8
$data[53,0] = 7
$data[53,1] = 7
$data[54,0] = 7
$data[54,1] = 7
$data[55,0] = 9
$data[55,1] = 9
$data[56,0] = 8
$data[56,1] = 8
$data[57,0] = 7
$data[57,1] = 7
$data[58,0] = 11
$data[58,1] = 11
$data[59,0] = 7
$data[59,1] = 7
$data[60,0] = 7
$data[60,1] = 7
$data[61,0] = 8
$data[61,1] = 8
$data[62,0] = 8
$data[62,1] = 8
$data[63,0] = 7
$data[63,1] = 7
$data[64,0] = 8
$data[64,1] = 8
$data[65,0] = 8
$data[65,1] = 8
$data[66,0] = 8
$data[66,1] = 8
$data[67,0] = 7
$data[67,1] = 7
$data[68,0] = 8
$data[68,1] = 8
$data[69,0] = 10
$data[69,1] = 10
$data[70,0] = 8
$data[70,1] = 8
$data[71,0] = 8
$data[71,1] = 8
$data[72,0] = 8
$data[72,1] = 8
$data[73,0] = 7
$data[73,1] = 7
$data[74,0] = 8
$data[74,1] = 8
$data[75,0] = 8
$data[75,1] = 8
$data[76,0] = 6
$data[76,1] = 6
$data[77,0] = 8
$data[77,1] = 8
$data[78,0] = 6
$data[78,1] = 6
$data[79,0] = 9
$data[79,1] = 9
$data[80,0] = 5
$data[80,1] = 5
$data[81,0] = 4
$data[81,1] = 4

$ws.Range("I2:J83").Value = $data
